# Auto-generated edit script applying the 2025-10-29 08:35 data refresh
$wb = $excel.ActiveWorkbook

# --- Metadata sheet: bump the "Last Updated" timestamp ---
$wsMeta = $wb.Worksheets.Item("Metadata")
$wsMeta.Range("A2").Value = "29 Oct 2025, 08:35 AM"

# --- Top Losers sheet: refreshed Latest/Weekly/Monthly %, re-sorted rows ---
$wsLosers = $wb.Worksheets.Item("Top Losers")
$wsLosers.Range("C2").Value = -9.319800000000001
$wsLosers.Range("D2").Value = -5.5246
$wsLosers.Range("E2").Value = 20.9934
$wsLosers.Range("C3").Value = -8.206300000000001
$wsLosers.Range("D3").Value = -9.955500000000001
$wsLosers.Range("E3").Value = -10.2321
$wsLosers.Range("C4").Value = -7.9006
$wsLosers.Range("D4").Value = 13.4399
$wsLosers.Range("E4").Value = 17.8523
$wsLosers.Range("C5").Value = -7.6698
$wsLosers.Range("D5").Value = -2.7601
$wsLosers.Range("E5").Value = 12.7574
$wsLosers.Range("C6").Value = -6.5228
$wsLosers.Range("D6").Value = 5.7921
$wsLosers.Range("E6").Value = 7.7813
$wsLosers.Range("C7").Value = -5.7681
$wsLosers.Range("D7").Value = -2.1511
$wsLosers.Range("E7").Value = 2.8599
$wsLosers.Range("C8").Value = -5.4482
$wsLosers.Range("D8").Value = -6.8581
$wsLosers.Range("E8").Value = 23.6844
$wsLosers.Range("C11").Value = -5.0058
$wsLosers.Range("D11").Value = -6.8383
$wsLosers.Range("E11").Value = -1.0019
$wsLosers.Range("C12").Value = -4.7668
$wsLosers.Range("D12").Value = 5.978
$wsLosers.Range("B13").Value = "HDFCAMC"
$wsLosers.Range("C13").Value = -4.4541
$wsLosers.Range("D13").Value = -2.6788
$wsLosers.Range("E13").Value = -2.4853
$wsLosers.Range("B14").Value = "CCCL"
$wsLosers.Range("C14").Value = -4.4296
$wsLosers.Range("D14").Value = -3.9882
$wsLosers.Range("E14").Value = -11.9368
$wsLosers.Range("C15").Value = -4.2305
$wsLosers.Range("D15").Value = -2.1142
$wsLosers.Range("E15").Value = 6.9792
$wsLosers.Range("C16").Value = -4.1048
$wsLosers.Range("D16").Value = -5.229
$wsLosers.Range("E16").Value = 9.7363
$wsLosers.Range("C17").Value = -4.0695
$wsLosers.Range("D17").Value = -11.7717
$wsLosers.Range("E17").Value = 16.7287
$wsLosers.Range("B18").Value = "VERANDA"
$wsLosers.Range("C18").Value = -3.9315
$wsLosers.Range("D18").Value = -4.2481
$wsLosers.Range("E18").Value = 11.7566
$wsLosers.Range("C19").Value = -3.8293
$wsLosers.Range("D19").Value = 6.025
$wsLosers.Range("E19").Value = -1.8685
$wsLosers.Range("B20").Value = "CAMS"
$wsLosers.Range("C20").Value = -3.8191
$wsLosers.Range("D20").Value = -1.2164
$wsLosers.Range("E20").Value = 1.9795
$wsLosers.Range("B21").Value = "CHENNPETRO"
$wsLosers.Range("C21").Value = -3.712
$wsLosers.Range("D21").Value = 5.1207
$wsLosers.Range("E21").Value = 7.1391
$wsLosers.Range("B22").Value = "STARHEALTH"
$wsLosers.Range("C22").Value = -3.6794
$wsLosers.Range("D22").Value = -2.3768
$wsLosers.Range("E22").Value = 6.648
$wsLosers.Range("B23").Value = "KALAMANDIR"
$wsLosers.Range("C23").Value = -3.6586
$wsLosers.Range("D23").Value = 3.0099
$wsLosers.Range("E23").Value = 27.5659
$wsLosers.Range("B24").Value = "PRIVISCL"
$wsLosers.Range("C24").Value = -3.6073
$wsLosers.Range("D24").Value = -3.0886
$wsLosers.Range("E24").Value = 18.5417
$wsLosers.Range("B25").Value = "BOSCHLTD"
$wsLosers.Range("C25").Value = -3.5807
$wsLosers.Range("D25").Value = -3.6931
$wsLosers.Range("E25").Value = -2.5823
$wsLosers.Range("B26").Value = "KHAICHEM"
$wsLosers.Range("C26").Value = -3.5409
$wsLosers.Range("D26").Value = -1.7041
$wsLosers.Range("E26").Value = 0.4825
$wsLosers.Range("B27").Value = "ABSLAMC"
$wsLosers.Range("C27").Value = -3.5004
$wsLosers.Range("D27").Value = -5.9054
$wsLosers.Range("E27").Value = -1.2571
$wsLosers.Range("B28").Value = "ENDURANCE"
$wsLosers.Range("C28").Value = -3.4618
$wsLosers.Range("D28").Value = -2.8209
$wsLosers.Range("E28").Value = 2.8958
$wsLosers.Range("B29").Value = "FABTECH"
$wsLosers.Range("C29").Value = -3.4572
$wsLosers.Range("D29").Value = 18.9369
$wsLosers.Range("E29").Value = "N/A"
$wsLosers.Range("B30").Value = "ANANDRATHI"
$wsLosers.Range("C30").Value = -3.3695
$wsLosers.Range("D30").Value = -1.1659
$wsLosers.Range("E30").Value = 8.8545
$wsLosers.Range("B31").Value = "SMSPHARMA"
$wsLosers.Range("C31").Value = -3.2801
$wsLosers.Range("D31").Value = -2.6303
$wsLosers.Range("E31").Value = 17.9923
$wsLosers.Range("B32").Value = "DIGITIDE"
$wsLosers.Range("C32").Value = -3.2703
$wsLosers.Range("D32").Value = 2.8163
$wsLosers.Range("E32").Value = 5.869
$wsLosers.Range("B33").Value = "INDIQUBE"
$wsLosers.Range("C33").Value = -3.268
$wsLosers.Range("D33").Value = -4.0347
$wsLosers.Range("E33").Value = -6.0088
$wsLosers.Range("B34").Value = "SHAREINDIA"
$wsLosers.Range("C34").Value = -3.2664
$wsLosers.Range("D34").Value = -0.8544
$wsLosers.Range("E34").Value = 56.035
$wsLosers.Range("B35").Value = "ATHERENERG"
$wsLosers.Range("C35").Value = -3.2525
$wsLosers.Range("D35").Value = 0.8637
$wsLosers.Range("E35").Value = 25.977
$wsLosers.Range("C36").Value = -3.2033
$wsLosers.Range("D36").Value = -4.8323
$wsLosers.Range("E36").Value = -7.5828
$wsLosers.Range("C37").Value = -3.1268
$wsLosers.Range("D37").Value = 2.0051
$wsLosers.Range("E37").Value = 1.0282
$wsLosers.Range("C38").Value = -3.0944
$wsLosers.Range("D38").Value = -1.3249
$wsLosers.Range("E38").Value = 6.3576
$wsLosers.Range("B39").Value = "RMDRIP"
$wsLosers.Range("C39").Value = -3.0679
$wsLosers.Range("D39").Value = -0.8073
$wsLosers.Range("E39").Value = 2.0567
$wsLosers.Range("B40").Value = "BHARATWIRE"
$wsLosers.Range("C40").Value = -3.0267
$wsLosers.Range("D40").Value = 23.4779
$wsLosers.Range("E40").Value = 24.5477
$wsLosers.Range("B41").Value = "DOLPHIN"
$wsLosers.Range("C41").Value = -2.995
$wsLosers.Range("D41").Value = -2.1521
$wsLosers.Range("E41").Value = -6.0936
$wsLosers.Range("C42").Value = -2.9699
$wsLosers.Range("D42").Value = 2.9393
$wsLosers.Range("E42").Value = 13.9851
$wsLosers.Range("B43").Value = "EDELWEISS"
$wsLosers.Range("C43").Value = -2.952
$wsLosers.Range("D43").Value = -3.7808
$wsLosers.Range("E43").Value = 8.0741
$wsLosers.Range("B44").Value = "CREST"
$wsLosers.Range("C44").Value = -2.9044
$wsLosers.Range("D44").Value = -5.5612
$wsLosers.Range("E44").Value = 3.0265
$wsLosers.Range("B45").Value = "SPARC"
$wsLosers.Range("C45").Value = -2.9037
$wsLosers.Range("D45").Value = 5.1229
$wsLosers.Range("E45").Value = 6.6245
$wsLosers.Range("B46").Value = "TIL"
$wsLosers.Range("C46").Value = -2.8939
$wsLosers.Range("D46").Value = -0.7558
$wsLosers.Range("E46").Value = -3.3136
$wsLosers.Range("B47").Value = "DRREDDY"
$wsLosers.Range("C47").Value = -2.8075
$wsLosers.Range("D47").Value = -2.3683
$wsLosers.Range("E47").Value = 2.4107
$wsLosers.Range("B48").Value = "KIRIINDUS"
$wsLosers.Range("C48").Value = -2.8066
$wsLosers.Range("D48").Value = -1.3498
$wsLosers.Range("E48").Value = 1.4695
$wsLosers.Range("C50").Value = -2.7248
$wsLosers.Range("D50").Value = -8.076599999999999
$wsLosers.Range("E50").Value = -8.0687
$wsLosers.Range("B51").Value = "SAMHI"
$wsLosers.Range("C51").Value = -2.6734
$wsLosers.Range("D51").Value = 1.6716
$wsLosers.Range("E51").Value = 2.6986
$wsLosers.Range("B52").Value = "NLCINDIA"
$wsLosers.Range("C52").Value = -2.6336
$wsLosers.Range("D52").Value = -4.1264
$wsLosers.Range("E52").Value = -11.24
$wsLosers.Range("B53").Value = "PRUDENT"
$wsLosers.Range("C53").Value = -2.5494
$wsLosers.Range("D53").Value = -2.935
$wsLosers.Range("E53").Value = 2.7302
$wsLosers.Range("B54").Value = "TTKPRESTIG"
$wsLosers.Range("C54").Value = -2.5414
$wsLosers.Range("D54").Value = 8.226100000000001
$wsLosers.Range("E54").Value = 9.8788
$wsLosers.Range("C55").Value = -2.4861
$wsLosers.Range("D55").Value = -3.5357
$wsLosers.Range("E55").Value = 1.0753
$wsLosers.Range("B56").Value = "PILANIINVS"
$wsLosers.Range("C56").Value = -2.4546
$wsLosers.Range("D56").Value = -0.7907
$wsLosers.Range("E56").Value = 4.267
$wsLosers.Range("B57").Value = "YATRA"
$wsLosers.Range("C57").Value = -2.4387
$wsLosers.Range("D57").Value = -2.2427
$wsLosers.Range("E57").Value = 8.0372
$wsLosers.Range("B58").Value = "NSIL"
$wsLosers.Range("C58").Value = -2.4088
$wsLosers.Range("D58").Value = -1.7646
$wsLosers.Range("E58").Value = 4.7431
$wsLosers.Range("B59").Value = "AYMSYNTEX"
$wsLosers.Range("C59").Value = -2.3918
$wsLosers.Range("D59").Value = 0.1563
$wsLosers.Range("E59").Value = -10.0208
$wsLosers.Range("B60").Value = "SUBROS"
$wsLosers.Range("C60").Value = -2.3713
$wsLosers.Range("D60").Value = 1.2248
$wsLosers.Range("E60").Value = 2.0455
$wsLosers.Range("B61").Value = "DAMCAPITAL"
$wsLosers.Range("C61").Value = -2.3557
$wsLosers.Range("D61").Value = -1.4973
$wsLosers.Range("E61").Value = 5.1146
$wsLosers.Range("B62").Value = "CUPID"
$wsLosers.Range("C62").Value = -2.3298
$wsLosers.Range("D62").Value = 1.9072
$wsLosers.Range("E62").Value = 11.3603
$wsLosers.Range("B63").Value = "ALLDIGI"
$wsLosers.Range("C63").Value = -2.3159
$wsLosers.Range("D63").Value = 0.0956
$wsLosers.Range("E63").Value = -5.0008
$wsLosers.Range("B64").Value = "SGFIN"
$wsLosers.Range("C64").Value = -2.213
$wsLosers.Range("D64").Value = 0.3261
$wsLosers.Range("E64").Value = 12.1582
$wsLosers.Range("B66").Value = "FCL"
$wsLosers.Range("C66").Value = -2.1612
$wsLosers.Range("D66").Value = -2.4325
$wsLosers.Range("E66").Value = 0.1684
$wsLosers.Range("B67").Value = "NUVAMA"
$wsLosers.Range("C67").Value = -2.1498
$wsLosers.Range("D67").Value = 1.4841
$wsLosers.Range("E67").Value = 15.5036
$wsLosers.Range("B68").Value = "ALICON"
$wsLosers.Range("C68").Value = -2.145
$wsLosers.Range("D68").Value = 5.7221
$wsLosers.Range("E68").Value = 10.9736
$wsLosers.Range("B69").Value = "BBOX"
$wsLosers.Range("C69").Value = -2.1376
$wsLosers.Range("D69").Value = -4.6405
$wsLosers.Range("E69").Value = 5.2887
$wsLosers.Range("B70").Value = "WEALTH"
$wsLosers.Range("C70").Value = -2.1352
$wsLosers.Range("D70").Value = -3.6938
$wsLosers.Range("E70").Value = -2.6549
$wsLosers.Range("B71").Value = "DEEDEV"
$wsLosers.Range("C71").Value = -2.1159
$wsLosers.Range("D71").Value = -6.4449
$wsLosers.Range("E71").Value = -7.2165
$wsLosers.Range("B72").Value = "KMEW"
$wsLosers.Range("C72").Value = -2.1113
$wsLosers.Range("D72").Value = 0.9812
$wsLosers.Range("E72").Value = -4.387
$wsLosers.Range("B73").Value = "CEATLTD"
$wsLosers.Range("C73").Value = -2.0911
$wsLosers.Range("D73").Value = -3.55
$wsLosers.Range("E73").Value = 18.044
$wsLosers.Range("B74").Value = "PENINLAND"
$wsLosers.Range("C74").Value = -2.0755
$wsLosers.Range("D74").Value = -0.1763
$wsLosers.Range("E74").Value = -1.4791
$wsLosers.Range("B75").Value = "VIKRAN"
$wsLosers.Range("C75").Value = -2.0326
$wsLosers.Range("D75").Value = 3.2747
$wsLosers.Range("E75").Value = 7.6826
$wsLosers.Range("B76").Value = "UTIAMC"
$wsLosers.Range("C76").Value = -2.0321
$wsLosers.Range("D76").Value = -4.6781
$wsLosers.Range("E76").Value = -2.0321

# --- 1 Month Performance sheet: refreshed % Change, re-sorted rows ---
$wsPerf = $wb.Worksheets.Item("1 Month Performance")
$wsPerf.Range("C7").Value = 66.6794
$wsPerf.Range("C8").Value = 64.4179
$wsPerf.Range("C10").Value = 52.5538
$wsPerf.Range("C12").Value = 45.3785
$wsPerf.Range("B13").Value = "TVSSRICHAK"
$wsPerf.Range("C13").Value = 40.3617
$wsPerf.Range("B14").Value = "MTARTECH"
$wsPerf.Range("C14").Value = 40.2667
$wsPerf.Range("C15").Value = 38.7665
$wsPerf.Range("B16").Value = "V2RETAIL"
$wsPerf.Range("C16").Value = 37.9127
$wsPerf.Range("B17").Value = "SEJALLTD"
$wsPerf.Range("C17").Value = 37.7023
$wsPerf.Range("C18").Value = 37.0114
$wsPerf.Range("B19").Value = "RAMAPHO"
$wsPerf.Range("C19").Value = 36.5309
$wsPerf.Range("B20").Value = "SHAREINDIA"
$wsPerf.Range("C20").Value = 36.4693
$wsPerf.Range("B21").Value = "NETWEB"
$wsPerf.Range("C21").Value = 36.4043
$wsPerf.Range("B22").Value = "TARACHAND"
$wsPerf.Range("C22").Value = 36.2607
$wsPerf.Range("C24").Value = 35.8
$wsPerf.Range("C29").Value = 32.4014
$wsPerf.Range("C34").Value = 27.83
$wsPerf.Range("B35").Value = "ARFIN"
$wsPerf.Range("C35").Value = 27.4033
$wsPerf.Range("B36").Value = "BHARATWIRE"
$wsPerf.Range("C36").Value = 27.1912
$wsPerf.Range("C37").Value = 27.1831
$wsPerf.Range("B39").Value = "HATSUN"
$wsPerf.Range("C39").Value = 26.4258
$wsPerf.Range("B40").Value = "AVALON"
$wsPerf.Range("C40").Value = 26.4142
$wsPerf.Range("B42").Value = "SCI"
$wsPerf.Range("C42").Value = 25.3374
$wsPerf.Range("B43").Value = "IFBIND"
$wsPerf.Range("C43").Value = 25.2709
$wsPerf.Range("C44").Value = 24.9094
$wsPerf.Range("C45").Value = 24.6277
$wsPerf.Range("C47").Value = 23.8061
$wsPerf.Range("C48").Value = 23.6586
$wsPerf.Range("B50").Value = "PRECWIRE"
$wsPerf.Range("C50").Value = 23.2681
$wsPerf.Range("B51").Value = "AUBANK"
$wsPerf.Range("C51").Value = 23.1988
$wsPerf.Range("B52").Value = "ETHOSLTD"
$wsPerf.Range("C52").Value = 23.0678
$wsPerf.Range("B53").Value = "ORBTEXP"
$wsPerf.Range("C53").Value = 22.9776
$wsPerf.Range("C55").Value = 22.2023
$wsPerf.Range("C56").Value = 21.9785
$wsPerf.Range("C57").Value = 21.8198
$wsPerf.Range("C61").Value = 20.2897
$wsPerf.Range("B64").Value = "ATL"
$wsPerf.Range("C64").Value = 19.7233
$wsPerf.Range("B65").Value = "BANKINDIA"
$wsPerf.Range("C65").Value = 19.6712
$wsPerf.Range("C66").Value = 19.5834
$wsPerf.Range("C67").Value = 19.4854
$wsPerf.Range("C68").Value = 19.4783
$wsPerf.Range("B69").Value = "CEATLTD"
$wsPerf.Range("C69").Value = 19.3991
$wsPerf.Range("B70").Value = "SUBROS"
$wsPerf.Range("C70").Value = 19.3491
$wsPerf.Range("C71").Value = 19.0575
$wsPerf.Range("C72").Value = 19.0526
$wsPerf.Range("C75").Value = 18.4848
$wsPerf.Range("B76").Value = "M&MFIN"
$wsPerf.Range("C76").Value = 18.4006
